$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row (44) appended at the bottom of the resale-number log.

# Columns A and D look like a date ("2024-01-10") and a zero-padded week
# number ("01"). Excel's normal auto-detection would turn these into a
# real date serial / a plain number, but the source data stores them as
# literal text (matching every other row in the sheet). Using a leading
# apostrophe forces Excel to keep them as text, then ClearFormats()
# removes the "quote prefix" cell style that the apostrophe trick adds,
# so the cell is left with the default (unstyled) formatting - just like
# the rest of the sheet.
$ws.Range("A44").Formula = "'2024-01-10"
$ws.Range("D44").Formula = "'01"
$ws.Range("A44").ClearFormats()
$ws.Range("D44").ClearFormats()

$ws.Range("B44").Value = "19:13:47"
$ws.Range("C44").Value = "Wednesday"

$ws.Range("E44").Value = 139484
$ws.Range("F44").Value = 142640
$ws.Range("G44").Value = 172153
$ws.Range("H44").Value = 148026
$ws.Range("I44").Value = -1
$ws.Range("J44").Value = 119236
$ws.Range("K44").Value = 224802
$ws.Range("L44").Value = 251428
$ws.Range("M44").Value = 185116
$ws.Range("N44").Value = 110468
$ws.Range("O44").Value = 40778
$ws.Range("P44").Value = 30870
$ws.Range("Q44").Value = 72785
$ws.Range("R44").Value = -1
$ws.Range("S44").Value = 42222
$ws.Range("T44").Value = -1
